{"js": "// Update the secsse model-comparison table: refreshed Log Likelihood / AIC\n// values for the CTD4, CTD3, MuHiSSE and CTD2 rows (plus the already-updated\n// MuSSE row), and the auto row-height hints that Word recalculated for the\n// CTD4, MuHiSSE and CTD2 rows as a result.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Column layout: 0=Model, 1=# params, 2=Log Likelihood, 3=AIC, 4=AIC weight\n// Replace just the number inside the cell so the run's existing formatting\n// (font, size, color, etc.) is preserved instead of being reset.\nasync function replaceInCell(rowIndex, colIndex, oldText, newText) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait replaceInCell(1, 2, \"192.15\", \"190.03\");     // CTD4    Log Likelihood\nawait replaceInCell(1, 3, \"-350.30\", \"-346.05\");   // CTD4    AIC\n\nawait replaceInCell(2, 2, \"177.15\", \"174.85\");     // CTD3    Log Likelihood\nawait replaceInCell(2, 3, \"-334.29\", \"-329.70\");   // CTD3    AIC\n\nawait replaceInCell(3, 2, \"139.77\", \"137.54\");     // MuHiSSE Log Likelihood\nawait replaceInCell(3, 3, \"-253.54\", \"-249.08\");   // MuHiSSE AIC\n\nawait replaceInCell(4, 2, \"123.96\", \"121.63\");     // CTD2    Log Likelihood\nawait replaceInCell(4, 3, \"-237.93\", \"-233.27\");   // CTD2    AIC\n\nawait replaceInCell(5, 2, \"-220.11\", \"-221.84\");   // MuSSE   Log Likelihood\nawait replaceInCell(5, 3, \"452.23\", \"455.67\");     // MuSSE   AIC\n\nawait context.sync();\n\n// Row heights (twips/20 = points): 638->635, 635 unchanged, 635->634, 637->634, 634 unchanged\nrows.items[1].preferredHeight = 635 / 20;  // body1 (CTD4)\nrows.items[3].preferredHeight = 634 / 20;  // body3 (MuHiSSE)\nrows.items[4].preferredHeight = 634 / 20;  // body4 (CTD2)\n\nawait context.sync();\n", "ps1": "# Update the secsse model-comparison table: refreshed Log Likelihood / AIC\n# values for the CTD4, CTD3, MuHiSSE and CTD2 rows (plus the already-updated\n# MuSSE row), and the auto row-height hints that Word recalculated for the\n# CTD4, MuHiSSE and CTD2 rows as a result.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Column layout: 1=Model, 2=# params, 3=Log Likelihood, 4=AIC, 5=AIC weight\n$t.Cell(2, 3).Range.Text = \"190.03\"   # CTD4    Log Likelihood\n$t.Cell(2, 4).Range.Text = \"-346.05\"  # CTD4    AIC\n\n$t.Cell(3, 3).Range.Text = \"174.85\"   # CTD3    Log Likelihood\n$t.Cell(3, 4).Range.Text = \"-329.70\"  # CTD3    AIC\n\n$t.Cell(4, 3).Range.Text = \"137.54\"   # MuHiSSE Log Likelihood\n$t.Cell(4, 4).Range.Text = \"-249.08\"  # MuHiSSE AIC\n\n$t.Cell(5, 3).Range.Text = \"121.63\"   # CTD2    Log Likelihood\n$t.Cell(5, 4).Range.Text = \"-233.27\"  # CTD2    AIC\n\n$t.Cell(6, 3).Range.Text = \"-221.84\"  # MuSSE   Log Likelihood\n$t.Cell(6, 4).Range.Text = \"455.67\"   # MuSSE   AIC\n\n# Row heights (twips/20 = points): 638->635, 635 unchanged, 635->634, 637->634, 634 unchanged\n$t.Rows.Item(2).Height = 635 / 20     # body1 (CTD4)\n$t.Rows.Item(4).Height = 634 / 20     # body3 (MuHiSSE)\n$t.Rows.Item(5).Height = 634 / 20     # body4 (CTD2)\n"}
